# Update specific numeric cell values on Sheet1 to reflect new algorithm
# results (commit: "Update Name of Algo").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value  = 6.839000000000001
$ws.Range("B3").Value  = 6.947
$ws.Range("E4").Value  = 13.038
$ws.Range("B5").Value  = 6.87
$ws.Range("E6").Value  = 12.754
$ws.Range("D7").Value  = -7.105
$ws.Range("A9").Value  = -21.393
$ws.Range("D9").Value  = -7.917999999999999
$ws.Range("E10").Value = 12.341
$ws.Range("B11").Value = 6.87
$ws.Range("E11").Value = 12.045
$ws.Range("B12").Value = 6.93
$ws.Range("A13").Value = -21.832
$ws.Range("A16").Value = -20.945
$ws.Range("A18").Value = -21.682
$ws.Range("A20").Value = -21.661
$ws.Range("B21").Value = 7.209999999999999
$ws.Range("D21").Value = -7.526999999999999
$ws.Range("E21").Value = 12.261
$ws.Range("E25").Value = 12.498

$wb.Save()
